# Laborator 04.04 - Jocul Minesweeper, folosind Parcurgerea in Adancime
# Mark week-7 attendance (column H) as present (TRUE) for the students
# that attended that week's lab.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(7, 8, 9, 14, 16, 17, 25, 28, 31, 32, 37, 41, 42, 45, 46, 47)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 8).Value = $true
}

# Reflect where the sheet was scrolled to / selected when the author saved.
$ws.Range("L25").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
